$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values per row (row index => column letter => value)
$data = @{
    2 = @{ E=3; G=203.7816646666667; H=611.344994; I=0.6667327591988204; J=0.6667327591988205; K=3; M=24.519512; N=73.558536; O=0.4736537296697991; P=0.4736537296697991; Q=4996.62697217431; R=44969.64274956879; S=0.3158004580875573; T=0.3158004580875574 }
    3 = @{ E=3; G=203.7816646666667; H=611.344994; I=0.6667327591988204; J=0.6667327591988205; K=3; M=14.70328633333333; N=44.109859; O=0.284029568377475; P=0.284029568377475; Q=2996.260165077316; R=26966.34148569585; S=0.1893718178183639; T=0.1893718178183639 }
    4 = @{ E=3; G=203.7816646666667; H=611.344994; I=0.6667327591988204; J=0.6667327591988205; K=3; M=12.543947; N=37.631841; O=0.2423167019527259; P=0.2423167019527259; Q=2556.226401150439; R=23006.03761035396; S=0.1615604832928991; T=0.1615604832928992 }
    5 = @{ E=3; G=63.14058933333333; H=189.421768; I=0.2065833519051582; J=0.2065833519051582; K=3; M=24.519512; N=73.558536; O=0.4736537296697991; P=0.4736537296697991; Q=1548.176437845739; R=13933.58794061165; S=0.09784897511756675; T=0.09784897511756677 }
    6 = @{ E=3; G=63.14058933333333; H=189.421768; I=0.2065833519051582; J=0.2065833519051582; K=3; M=14.70328633333333; N=44.109859; O=0.284029568377475; P=0.284029568377475; Q=928.3741642234123; R=8355.367478010712; S=0.05867578027559409; T=0.0586757802755941 }
    7 = @{ E=3; G=63.14058933333333; H=189.421768; I=0.2065833519051582; J=0.2065833519051582; K=3; M=12.543947; N=37.631841; O=0.2423167019527259; P=0.2423167019527259; Q=792.0322061460987; R=7128.289855314888; S=0.05005859651199731; T=0.05005859651199732 }
    8 = @{ E=3; G=38.719942; H=116.159826; I=0.1266838888960214; J=0.1266838888960214; K=3; M=24.519512; N=73.558536; O=0.4736537296697991; P=0.4736537296697991; Q=949.394082508304; R=8544.546742574736; S=0.06000429646467498; T=0.06000429646467499 }
    9 = @{ E=3; G=38.719942; H=116.159826; I=0.1266838888960214; J=0.1266838888960214; K=3; M=14.70328633333333; N=44.109859; O=0.284029568377475; P=0.284029568377475; Q=569.3103940360593; R=5123.793546324534; S=0.03598197028351695; T=0.03598197028351695 }
    10 = @{ E=3; G=38.719942; H=116.159826; I=0.1266838888960214; J=0.1266838888960214; K=3; M=12.543947; N=37.631841; O=0.2423167019527259; P=0.2423167019527259; Q=485.700900291074; R=4371.308102619666; S=0.03069762214782946; T=0.03069762214782947 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
